$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing the existing rows 21-22 down to 22-23.
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new weekly price observation.
$ws.Cells.Item(21, 1).Value = 6
$ws.Cells.Item(21, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(21, 3).Value = 'Metropolitana'
$ws.Cells.Item(21, 4).Value = Get-Date -Year 2021 -Month 11 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = 'Fruta'
$ws.Cells.Item(21, 7).Value = 100102
$ws.Cells.Item(21, 8).Value = 'Cítricos'
$ws.Cells.Item(21, 9).Value = 100102006
$ws.Cells.Item(21, 10).Value = 'Pomelo'
$ws.Cells.Item(21, 11).Value = 'Start Ruby'
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 24
$ws.Cells.Item(21, 14).Value = 140000
$ws.Cells.Item(21, 15).Value = 150000
$ws.Cells.Item(21, 16).Value = 145000
$ws.Cells.Item(21, 17).Value = '$/bins (350 kilos)'
$ws.Cells.Item(21, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 19).Value = 414
$ws.Cells.Item(21, 20).Value = 350
